# ------------------------------------------------------------------
# feat: add 2022-Q4 data
#
# 1) Insert a brand-new worksheet named "2022-Q4" right after "总计",
#    pushing the existing quarter sheets (2022-Q3 .. 2021-Q2) one slot
#    to the right. Populate it with the fund holding detail rows
#    (fund code / name / scale / position / ratio / market value are
#    stored as plain text, matching how the other quarter sheets do
#    it; only the row index and the rank column are real numbers),
#    reusing the "2022-Q3" sheet's header/index-column formatting so
#    the new sheet looks consistent with its siblings.
# 2) Update the "总计" (summary) sheet: insert a new data row for
#    2022-Q4 at the top of the table and shift every other quarter's
#    row down by one, appending the final 2021-Q2 row at the bottom.
# ------------------------------------------------------------------

# Write a value as plain text (no leading-zero / trailing-zero loss)
# while leaving the cell on the default "Normal" style afterwards, so
# no stray numeric-format style gets attached to the cell.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")

# --- 1) create + populate the new "2022-Q4" sheet, placed right after "总计" ---
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# NOTE: re-fetch "2022-Q3" after Add() -- the handle captured before the
# insert can otherwise resolve to a stale sheet position and silently
# copy nothing.
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Range("B1:H1").Copy($q4.Range("B1:H1"))   # header formatting (s="2")
$q3.Range("A2").Copy($q4.Range("A2:A13"))     # index-column formatting (s="2")

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$q4data = @(
    @(0,  "159865", "国泰中证畜牧养殖ETF",                    "34.10", "99.68", "2.34", "0.7979", 10),
    @(1,  "501201", "红土创新科技创新 3 年封闭",                "3.82",  "92.71", "4.57", "0.1746", 5),
    @(2,  "014179", "中银证券远见价值混合A",                    "1.59",  "85.40", "8.92", "0.1418", 1),
    @(3,  "516670", "招商中证畜牧养殖ETF",                      "6.00",  "99.30", "2.29", "0.1374", 10),
    @(4,  "159867", "鹏华中证畜牧养殖ETF",                      "5.22",  "97.53", "2.25", "0.1174", 10),
    @(5,  "003980", "中银证券瑞益灵活配置混合A",                "0.58",  "88.28", "8.90", "0.0516", 1),
    @(6,  "005571", "中银证券新能源灵活配置混合A",              "0.55",  "89.79", "7.84", "0.0431", 2),
    @(7,  "168401", "红土创新转型精选灵活配置混合（LOF）",       "0.82",  "93.08", "4.67", "0.0383", 5),
    @(8,  "516760", "平安中证畜牧养殖ETF",                      "1.51",  "98.80", "2.28", "0.0344", 10),
    @(9,  "005572", "中银证券新能源灵活配置混合C",              "0.26",  "89.79", "7.84", "0.0204", 2),
    @(10, "003981", "中银证券瑞益灵活配置混合C",                "0.21",  "88.28", "8.90", "0.0187", 1),
    @(11, "014180", "中银证券远见价值混合C",                    "0.13",  "85.40", "8.92", "0.0116", 1)
)

$r = 2
foreach ($row in $q4data) {
    $q4.Cells.Item($r, 1).Value = $row[0]                      # A: index (number)
    Set-TextValue $q4.Cells.Item($r, 2) $row[1]                # B: 基金代码 (text)
    Set-TextValue $q4.Cells.Item($r, 3) $row[2]                # C: 基金名称 (text)
    Set-TextValue $q4.Cells.Item($r, 4) $row[3]                # D: 基金规模 (text)
    Set-TextValue $q4.Cells.Item($r, 5) $row[4]                # E: 股票总仓位 (text)
    Set-TextValue $q4.Cells.Item($r, 6) $row[5]                # F: 仓位占比 (text)
    Set-TextValue $q4.Cells.Item($r, 7) $row[6]                # G: 持有市值(亿元) (text)
    $q4.Cells.Item($r, 8).Value = $row[7]                      # H: 仓位排名 (number)
    $r = $r + 1
}

# --- 2) update the "总计" summary sheet: shift rows down, insert 2022-Q4 on top ---

# Row 8 is brand new (table used to stop at row 7) -- clone column A's
# index-cell formatting from the previous last row (row 7) before the
# row-by-row value fill below overwrites row 7 with its new (shifted)
# contents, so the new row keeps the same "index" style as its peers.
$summary.Cells.Item(7, 1).Copy($summary.Cells.Item(8, 1))

$summaryData = @(
    @(0, "2022-Q4", 12, 1.59),
    @(1, "2022-Q3", 7,  0.26),
    @(2, "2022-Q2", 6,  0.26),
    @(3, "2022-Q1", 1,  0.01),
    @(4, "2021-Q4", 3,  0.11),
    @(5, "2021-Q3", 7,  0.49),
    @(6, "2021-Q2", 2,  0.01)
)

$r = 2
foreach ($row in $summaryData) {
    $summary.Cells.Item($r, 1).Value = $row[0]
    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$summary.Activate()
